# Update "相談件数" (consultation count) sheet with the next day's data row.
# The existing footer/note row (currently row 99, a merged-style text cell
# referencing the "※4/8..." note) is pushed down to row 100, and a new data
# row is written in its place at row 99.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("相談件数")

# Capture footer row 99 (text cell) before we overwrite it, then move it to row 100.
$footerValue = $ws.Range("B99").Value

$ws.Range("B100").Value = $footerValue

# Clear the old footer cell location so it no longer carries the text.
$ws.Range("B99").ClearContents()

# New data row 99.
$ws.Range("A99").Value = 43954
$ws.Range("B99").Value = 308
$ws.Range("C99").Value = 33036
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 6958

# Keep number formats consistent with the row above it (row 98).
$ws.Range("A99").NumberFormat = $ws.Range("A98").NumberFormat
$ws.Range("B99:C99").NumberFormat = $ws.Range("B98:C98").NumberFormat
$ws.Range("D99:E99").NumberFormat = $ws.Range("D98:E98").NumberFormat

# Update the print area to include the newly added row.
$ws.PageSetup.PrintArea = "`$A`$1:`$E`$101"

# Update frozen-pane anchor and active selection to match the new scroll position.
$ws.Activate()
$ws.Range("D80").Select()
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("B2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("D99").Select()
